# Apply the update described by the diff.
# Rows 11 and 13 effectively trade most of their descriptive data (species
# info, substrate info, times, public comment) while getting brand-new
# Id/Taxonsorteringsordning values. Rows 12 and 14 only get a new
# Taxonsorteringsordning (column B) value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: becomes the "Ulltickeporing / Skeletocutis brevispora" record ---
$ws.Range("A11").Value = 112243622
$ws.Range("B11").Value = 89893
$ws.Range("D11").Value = "VU"
$ws.Range("E11").Value = 2062
$ws.Range("F11").Value = "Ulltickeporing"
$ws.Range("G11").Value = "Skeletocutis brevispora"
$ws.Range("H11").Value = "Niemelä"
$ws.Range("Z11").Value = "10:21"
$ws.Range("AB11").Value = "10:21"
$ws.Range("AC11").Value = "På ytmurken granlåga med delvis avfallande bark och insektsgnag i veden från tiden då granen stod upp."
$ws.Range("AJ11").Value = "ullticka"
$ws.Range("AK11").Value = "Phellinidium ferrugineofuscum"
$ws.Range("AO11").Value = "Phellinidium ferrugineofuscum"

# --- Row 12: only Taxonsorteringsordning changes ---
$ws.Range("B12").Value = 77685

# --- Row 13: becomes the "Ullticka / Phellinidium ferrugineofuscum" record ---
$ws.Range("A13").Value = 112243230
$ws.Range("B13").Value = 89553
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 1202
$ws.Range("F13").Value = "Ullticka"
$ws.Range("G13").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H13").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Z13").Value = "10:16"
$ws.Range("AB13").Value = "10:16"
$ws.Range("AC13").ClearContents()
$ws.Range("AJ13").Value = "gran"
$ws.Range("AK13").Value = "Picea abies"
$ws.Range("AO13").Value = "Picea abies"

# --- Row 14: only Taxonsorteringsordning changes ---
$ws.Range("B14").Value = 88967
